# Update the NATMI LR-pairs sheet with recalculated TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9167996666666666
$ws.Range("H2").Value = 2.750399
$ws.Range("I2").Value = 0.2833456974325495
$ws.Range("J2").Value = 0.2833456974325495
$ws.Range("M2").Value = 2.709957333333333
$ws.Range("N2").Value = 8.129871999999999
$ws.Range("O2").Value = 0.03910908881921663
$ws.Range("P2").Value = 0.03910908881921663
$ws.Range("Q2").Value = 2.484487979880889
$ws.Range("R2").Value = 22.36039181892799
$ws.Range("S2").Value = 0.01108139204743246
$ws.Range("T2").Value = 0.01108139204743246

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9167996666666666
$ws.Range("H3").Value = 2.750399
$ws.Range("I3").Value = 0.2833456974325495
$ws.Range("J3").Value = 0.2833456974325495
$ws.Range("O3").Value = 0.868442511192471
$ws.Range("P3").Value = 0.868442511192471
$ws.Range("Q3").Value = 55.16965609321721
$ws.Range("R3").Value = 496.526904838955
$ws.Range("S3").Value = 0.2460694490139054
$ws.Range("T3").Value = 0.2460694490139054

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9167996666666666
$ws.Range("H4").Value = 2.750399
$ws.Range("I4").Value = 0.2833456974325495
$ws.Range("J4").Value = 0.2833456974325495
$ws.Range("O4").Value = 0.09244839998831231
$ws.Range("P4").Value = 0.09244839998831231
$ws.Range("Q4").Value = 5.872981075880332
$ws.Range("R4").Value = 52.856829682923
$ws.Range("S4").Value = 0.02619485637121165
$ws.Range("T4").Value = 0.02619485637121165

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.2271242616180895
$ws.Range("J5").Value = 0.2271242616180895
$ws.Range("M5").Value = 2.709957333333333
$ws.Range("N5").Value = 8.129871999999999
$ws.Range("O5").Value = 0.03910908881921663
$ws.Range("P5").Value = 0.03910908881921663
$ws.Range("Q5").Value = 1.991516028097778
$ws.Range("R5").Value = 17.92364425288
$ws.Range("S5").Value = 0.008882622920620857
$ws.Range("T5").Value = 0.008882622920620857

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.2271242616180895
$ws.Range("J6").Value = 0.2271242616180895
$ws.Range("O6").Value = 0.868442511192471
$ws.Range("P6").Value = 0.868442511192471
$ws.Range("S6").Value = 0.1972443641123494
$ws.Range("T6").Value = 0.1972443641123494

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.2271242616180895
$ws.Range("J7").Value = 0.2271242616180895
$ws.Range("O7").Value = 0.09244839998831231
$ws.Range("P7").Value = 0.09244839998831231
$ws.Range("S7").Value = 0.02099727458511923
$ws.Range("T7").Value = 0.02099727458511923

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 1.583934333333334
$ws.Range("H8").Value = 4.751803000000001
$ws.Range("I8").Value = 0.4895300409493609
$ws.Range("J8").Value = 0.4895300409493609
$ws.Range("M8").Value = 2.709957333333333
$ws.Range("N8").Value = 8.129871999999999
$ws.Range("O8").Value = 0.03910908881921663
$ws.Range("P8").Value = 0.03910908881921663
$ws.Range("Q8").Value = 4.292394462135111
$ws.Range("R8").Value = 38.631550159216
$ws.Range("S8").Value = 0.01914507385116331
$ws.Range("T8").Value = 0.01914507385116331

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 1.583934333333334
$ws.Range("H9").Value = 4.751803000000001
$ws.Range("I9").Value = 0.4895300409493609
$ws.Range("J9").Value = 0.4895300409493609
$ws.Range("O9").Value = 0.868442511192471
$ws.Range("P9").Value = 0.868442511192471
$ws.Range("Q9").Value = 95.31538417979279
$ws.Range("R9").Value = 857.8384576181351
$ws.Range("S9").Value = 0.4251286980662161
$ws.Range("T9").Value = 0.4251286980662161

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 1.583934333333334
$ws.Range("H10").Value = 4.751803000000001
$ws.Range("I10").Value = 0.4895300409493609
$ws.Range("J10").Value = 0.4895300409493609
$ws.Range("O10").Value = 0.09244839998831231
$ws.Range("P10").Value = 0.09244839998831231
$ws.Range("R10").Value = 91.31956558223102
$ws.Range("S10").Value = 0.04525626903198142
$ws.Range("T10").Value = 0.04525626903198142
